# Refresh the cryptocurrency price/volume snapshot ("Updated cryptos list ... with GitHub Actions").
# Only the Price (D) and Volume(1h) (E) columns move for most rows; two rows
# (44 & 45) also swap their Coin name + Link because ImmutableX and OKB traded
# places in the ranking.
#
# Values that look numeric ("572.54", "1.01", ...) are written with a leading
# apostrophe so Excel stores them as text (matching the sheet's existing
# convention of keeping the Price column as text, e.g. "69.409.96" which is
# not even a valid number) instead of silently coercing them to doubles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2: Bitcoin
    $ws.Range("D2").Value = '69.537.11'
    $ws.Range("E2").Value = '  -1.56%  '

    # Row 3: Ethereum
    $ws.Range("D3").Value = '2.503.52'
    $ws.Range("E3").Value = '  -2.39%  '

    # Row 4: TetherUSD
    $ws.Range("E4").Value = '  +0.03%  '

    # Row 5: BNB
    $ws.Range("D5").Value = '''572.54'
    $ws.Range("E5").Value = '  -1.22%  '

    # Row 6: Solana
    $ws.Range("D6").Value = '''165.04'
    $ws.Range("E6").Value = '  -3.42%  '

    # Row 7: USDC
    $ws.Range("E7").Value = '  +0.06%  '

    # Row 8: XRP
    $ws.Range("D8").Value = '''0.517'
    $ws.Range("E8").Value = '  +1.05%  '

    # Row 9: LidoStakedEther
    $ws.Range("D9").Value = '2.500.50'
    $ws.Range("E9").Value = '  -2.49%  '

    # Row 10: Dogecoin
    $ws.Range("D10").Value = '''0.160'
    $ws.Range("E10").Value = '  -3.65%  '

    # Row 11: TRON
    $ws.Range("E11").Value = '  -1.18%  '

    # Row 12: Cardano
    $ws.Range("D12").Value = '''0.353'
    $ws.Range("E12").Value = '  +1.98%  '

    # Row 13: Toncoin
    $ws.Range("D13").Value = '''4.92'
    $ws.Range("E13").Value = '  +1.12%  '

    # Row 14: WrappedliquidstakedEther2.0
    $ws.Range("D14").Value = '2.958.97'
    $ws.Range("E14").Value = '  -2.46%  '

    # Row 15: WrappedBTC
    $ws.Range("D15").Value = '69.369.74'
    $ws.Range("E15").Value = '  -1.56%  '

    # Row 16: ShibaInu
    $ws.Range("D16").Value = '''0.0000176'
    $ws.Range("E16").Value = '  -3.92%  '

    # Row 17: Avalanche
    $ws.Range("D17").Value = '''24.95'
    $ws.Range("E17").Value = '  -1.32%  '

    # Row 18: WrappedEther
    $ws.Range("D18").Value = '2.516.90'
    $ws.Range("E18").Value = '  -1.92%  '

    # Row 19: Uniswap
    $ws.Range("D19").Value = '''7.81'
    $ws.Range("E19").Value = '  +4.70%  '

    # Row 20: Chainlink
    $ws.Range("D20").Value = '''11.34'
    $ws.Range("E20").Value = '  -3.74%  '

    # Row 21: BitcoinCash
    $ws.Range("D21").Value = '''347.41'
    $ws.Range("E21").Value = '  -4.45%  '

    # Row 22: Polkadot
    $ws.Range("D22").Value = '''3.90'
    $ws.Range("E22").Value = '  -2.64%  '

    # Row 23: SuiNetwork
    $ws.Range("E23").Value = '  -2.68%  '

    # Row 24: Dai
    $ws.Range("E24").Value = '  -0.08%  '

    # Row 25: Litecoin
    $ws.Range("D25").Value = '''70.04'
    $ws.Range("E25").Value = '  -0.17%  '

    # Row 26: NEARProtocol
    $ws.Range("D26").Value = '''3.97'
    $ws.Range("E26").Value = '  -3.99%  '

    # Row 27: WrappedeETH
    $ws.Range("D27").Value = '2.694.36'

    # Row 28: Aptos
    $ws.Range("D28").Value = '''8.77'
    $ws.Range("E28").Value = '  -6.71%  '

    # Row 29: Binance-PegBSC-USD
    $ws.Range("D29").Value = '''1.01'
    $ws.Range("E29").Value = '  +1.48%  '

    # Row 30: PEPE
    $ws.Range("D30").Value = '0.0₃0897'
    $ws.Range("E30").Value = '  -3.77%  '

    # Row 31: InternetComputer(DFINITY)
    $ws.Range("D31").Value = '''7.85'
    $ws.Range("E31").Value = '  -1.07%  '

    # Row 32: Bittensor
    $ws.Range("D32").Value = '''461.73'
    $ws.Range("E32").Value = '  -5.21%  '

    # Row 33: Fetch.AI
    $ws.Range("E33").Value = '  -5.02%  '

    # Row 34: PancakeSwap
    $ws.Range("D34").Value = '''1.73'
    $ws.Range("E34").Value = '  -2.26%  '

    # Row 35: FirstDigitalUSD
    $ws.Range("E35").Value = '  +0.11%  '

    # Row 36: Kaspa
    $ws.Range("E36").Value = '  +3.04%  '

    # Row 37: Monero
    $ws.Range("D37").Value = '''156.27'
    $ws.Range("E37").Value = '  -0.44%  '

    # Row 38: WhiteBITCoin
    $ws.Range("D38").Value = '''19.04'
    $ws.Range("E38").Value = '  +1.02%  '

    # Row 39: EthereumClassic
    $ws.Range("D39").Value = '''18.54'
    $ws.Range("E39").Value = '  -1.34%  '

    # Row 41: RenderToken
    $ws.Range("D41").Value = '''4.74'
    $ws.Range("E41").Value = '  -1.21%  '

    # Row 42: PolygonEcosystemToken
    $ws.Range("D42").Value = '''0.316'
    $ws.Range("E42").Value = '  -1.84%  '

    # Row 43: Stacks
    $ws.Range("D43").Value = '''1.60'
    $ws.Range("E43").Value = '  -5.76%  '

    # Row 44: OKB -> ImmutableX
    $ws.Range("B44").Value = 'ImmutableX'
    $ws.Range("C44").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    $ws.Range("D44").Value = '''1.15'
    $ws.Range("E44").Value = '  -14.11%  '

    # Row 45: ImmutableX -> OKB
    $ws.Range("B45").Value = 'OKB'
    $ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    $ws.Range("D45").Value = '''38.25'
    $ws.Range("E45").Value = '  -0.89%  '

    # Row 46: dogwifhat
    $ws.Range("D46").Value = '''2.29'
    $ws.Range("E46").Value = '  -8.48%  '

    # Row 47: Aave
    $ws.Range("D47").Value = '''142.25'
    $ws.Range("E47").Value = '  -2.57%  '

    # Row 48: ARBITRUM
    $ws.Range("D48").Value = '''0.524'
    $ws.Range("E48").Value = '  -1.95%  '

    # Row 49: Filecoin
    $ws.Range("D49").Value = '''3.47'
    $ws.Range("E49").Value = '  -3.15%  '

    # Row 50: Optimism
    $ws.Range("D50").Value = '''1.58'
    $ws.Range("E50").Value = '  -4.61%  '

    # Row 51: Cronos
    $ws.Range("D51").Value = '''0.0726'
    $ws.Range("E51").Value = '  -1.47%  '

